# Restore revision: update the "From" value for the RET1 rule (row 10)
# in the Rules sheet from 18 to 1.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

$ws.Range("C10").Value = 1
